$d = $word.ActiveDocument

function Replace-Once($old, $new) {
  $r = $d.Content
  $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
  if ($found) {
    $r.Text = $new
    return $true
  } else {
    return $false
  }
}

function Replace-All($old, $new) {
  $n = 0
  while (Replace-Once $old $new) {
    $n++
    if ($n -gt 50) { break }
  }
  if ($n -eq 0) {
    Write-Host "WARNING: no match found for:" $old
  }
  return $n
}

Replace-All "Wanahisabati wanaocheza:" "The playful mathematicians:" | Out-Null
Replace-All "** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino" "** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino" | Out-Null
Replace-All "[Muziki]" "[Music]" | Out-Null
Replace-All "kuna wanahisabati wawili, tupige simu" "there are two mathematicians, let's call" | Out-Null
Replace-All "Fil na Mike wanaokutana" "them Fil and Mike who meet each other" | Out-Null
Replace-All "tena baada ya muda mrefu. Baada ya baadhi" "again after a long time. After some" | Out-Null
Replace-All "kuzungumza, Phil anasema ana watoto watatu, basi" "chatting, Phil says he has three children, then" | Out-Null
Replace-All "Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil," "Mike, astonished, asks: 'How old are they?' Fil," | Out-Null
Replace-All "kuwa mwanahisabati mchezaji, anajibu" "being a playful mathematician, answers" | Out-Null
Replace-All "'Wewe niambie! Nitakupa kidokezo: ikiwa wewe" "'You tell me! I'll give you a hint: if you" | Out-Null
Replace-All "zidisheni enzi tatu pamoja ninyi" "multiply the three ages together you" | Out-Null
Replace-All "pata 36.' Mike huchukua wakati mwingine kufikiria" "get 36.' Mike takes sometimes to think" | Out-Null
Replace-All "na kusema: 'Samahani Fil, lakini nahitaji" "and says: 'I'm sorry Fil, but I do need" | Out-Null
Replace-All "kidokezo kingine. Kwa hivyo Fil anamwambia Mike:" "another hint. So Fil tells Mike:" | Out-Null
Replace-All "'Ndiyo, hakika, hapa ni: kama alikuwa na hadi" "'Yes, sure, here it is: if you had up to" | Out-Null
Replace-All "miaka mitatu unapata idadi ya hesabu" "three ages you get the number of math" | Out-Null
Replace-All "karatasi tunachapisha pamoja. Je, unaikumbuka?'" "papers we publish together. Do you remember it?'" | Out-Null
Replace-All "'Ndio nakumbuka wangapi, lakini bado" "'Yes I do remember How many, but still" | Out-Null
Replace-All "Sina taarifa za kutosha! nahitaji" "I do not have enough information! I need" | Out-Null
Replace-All "angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo" "at least one more.' Fil says: 'Yes don't" | Out-Null
Replace-All "wasiwasi lakini hii ni ya mwisho:" "worry but this is the last one:" | Out-Null
Replace-All "Mdogo ana macho ya blues.' Na" "The youngest one has blues eyes.' And" | Out-Null
Replace-All "ghafla Mike anapata jibu. Wewe" "suddenly Mike gets the answer. You" | Out-Null
Replace-All "sikia mazungumzo lakini hujui" "hear the conversation but you don't know" | Out-Null
Replace-All "ni karatasi ngapi walichapisha pamoja." "how many papers they published together." | Out-Null
Replace-All "Hata hivyo, unataka kujua umri wa" "However, you do want to know the ages of" | Out-Null
Replace-All "watoto watatu. Je, unaweza kuwahesabu" "the three children. Can you figure them" | Out-Null
Replace-All "nje?" "out?" | Out-Null
